{"js": "// Replace each \"old problem=\" text with its \"new problem=\" text inside the\n// document's multiplication-practice table. Every old value in this table\n// is unique, so an exact (non-wildcard, case-sensitive, whole-match) text\n// search safely targets the single cell that needs updating.\nconst replacements = [\n    [\"859\u00d77=\", \"325\u00d79=\"],\n    [\"528\u00d75=\", \"208\u00d78=\"],\n    [\"551\u00d75=\", \"867\u00d76=\"],\n    [\"458\u00d72=\", \"967\u00d74=\"],\n    [\"818\u00d75=\", \"545\u00d72=\"],\n    [\"625\u00d73=\", \"158\u00d76=\"],\n    [\"476\u00d73=\", \"726\u00d77=\"],\n    [\"345\u00d73=\", \"276\u00d74=\"],\n    [\"827\u00d78=\", \"392\u00d74=\"],\n    [\"872\u00d73=\", \"268\u00d75=\"],\n    [\"812\u00d73=\", \"677\u00d77=\"],\n    [\"365\u00d78=\", \"602\u00d75=\"],\n    [\"885\u00d77=\", \"142\u00d75=\"],\n    [\"481\u00d78=\", \"652\u00d72=\"],\n    [\"683\u00d77=\", \"815\u00d74=\"],\n    [\"102\u00d75=\", \"517\u00d74=\"],\n    [\"659\u00d72=\", \"535\u00d72=\"],\n    [\"988\u00d79=\", \"494\u00d73=\"],\n    [\"624\u00d73=\", \"143\u00d78=\"],\n    [\"665\u00d73=\", \"750\u00d76=\"],\n    [\"154\u00d79=\", \"346\u00d72=\"],\n    [\"838\u00d74=\", \"536\u00d73=\"],\n    [\"365\u00d79=\", \"380\u00d74=\"],\n    [\"465\u00d72=\", \"345\u00d77=\"],\n    [\"470\u00d77=\", \"516\u00d73=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n    const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n    found.load(\"items\");\n    await context.sync();\n\n    for (const range of found.items) {\n        range.insertText(newText, \"Replace\");\n    }\n    await context.sync();\n}\n", "ps1": "# Update the multiplication-practice table: each old \"NNN\u00d7N=\" problem text\n# is replaced by its new \"NNN\u00d7N=\" text. Every old value is unique across the\n# document, so a plain (non-wildcard) Find/Replace-All per pair safely hits\n# exactly the one cell that needs updating.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"859\u00d77=\"; New = \"325\u00d79=\" },\n    @{ Old = \"528\u00d75=\"; New = \"208\u00d78=\" },\n    @{ Old = \"551\u00d75=\"; New = \"867\u00d76=\" },\n    @{ Old = \"458\u00d72=\"; New = \"967\u00d74=\" },\n    @{ Old = \"818\u00d75=\"; New = \"545\u00d72=\" },\n    @{ Old = \"625\u00d73=\"; New = \"158\u00d76=\" },\n    @{ Old = \"476\u00d73=\"; New = \"726\u00d77=\" },\n    @{ Old = \"345\u00d73=\"; New = \"276\u00d74=\" },\n    @{ Old = \"827\u00d78=\"; New = \"392\u00d74=\" },\n    @{ Old = \"872\u00d73=\"; New = \"268\u00d75=\" },\n    @{ Old = \"812\u00d73=\"; New = \"677\u00d77=\" },\n    @{ Old = \"365\u00d78=\"; New = \"602\u00d75=\" },\n    @{ Old = \"885\u00d77=\"; New = \"142\u00d75=\" },\n    @{ Old = \"481\u00d78=\"; New = \"652\u00d72=\" },\n    @{ Old = \"683\u00d77=\"; New = \"815\u00d74=\" },\n    @{ Old = \"102\u00d75=\"; New = \"517\u00d74=\" },\n    @{ Old = \"659\u00d72=\"; New = \"535\u00d72=\" },\n    @{ Old = \"988\u00d79=\"; New = \"494\u00d73=\" },\n    @{ Old = \"624\u00d73=\"; New = \"143\u00d78=\" },\n    @{ Old = \"665\u00d73=\"; New = \"750\u00d76=\" },\n    @{ Old = \"154\u00d79=\"; New = \"346\u00d72=\" },\n    @{ Old = \"838\u00d74=\"; New = \"536\u00d73=\" },\n    @{ Old = \"365\u00d79=\"; New = \"380\u00d74=\" },\n    @{ Old = \"465\u00d72=\"; New = \"345\u00d77=\" },\n    @{ Old = \"470\u00d77=\"; New = \"516\u00d73=\" }\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.Text = $pair.New\n    $find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $pair.New, 2) | Out-Null\n}\n"}
